$d = $word.ActiveDocument

# Replace the first occurrence of $old found inside $range with $new.
# Uses Find to locate the precise span, then InsertAfter + a separate
# Delete of the original span (rather than a plain Find/Replace or a
# direct Range.Text assignment) because this interpreter re-anchors
# bookmark/comment markers that sit at a run boundary when a Replace
# or Text-assignment touches that exact position; inserting the new
# text first and deleting the old text afterwards keeps such markers
# (e.g. commentRangeStart) anchored where they belong.
function Find-Replace($range, $old, $new) {
    $r = $range.Duplicate()
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $old"
        return
    }
    $s = $r.Start
    $e = $r.End
    $r.InsertAfter($new)
    $toDelete = $word.ActiveDocument.Range($s, $e)
    $toDelete.Delete()
}

# Paragraph 3: standalone "English" heading (style P68B1DB1-Normal2).
# NB: leave the hyperlinked "English" in paragraph 1 untouched.
Find-Replace $d.Paragraphs.Item(3).Range "English" "Inglês"

# Paragraph 14: "Thank you for submitting your documents" heading.
Find-Replace $d.Paragraphs.Item(14).Range `
    "Thank you for submitting your documents" `
    "Obrigado por enviar os seus documentos"

# Paragraph 16: "Hi [PARTNER NAME], "
Find-Replace $d.Paragraphs.Item(16).Range "Hi " "Olá "
Find-Replace $d.Paragraphs.Item(16).Range "[PARTNER NAME]" "[NOME DO PARCEIRO]"

# Paragraph 18: intro + event name + arrangements sentence.
$para18 = $d.Paragraphs.Item(18).Range
Find-Replace $para18 `
    "Thank you for providing us with your documents for the upcoming " `
    "Obrigado por nos enviar os seus documentos para a "
Find-Replace $para18 "[EVENT NAME]" "[NOME DO EVENTO]"
Find-Replace $para18 `
    ". Based on the information you’ve given us, we’ll make the necessary arrangements, including accommodation and transportation." `
    ". Com base nas informações que nos enviou, iremos tomar as providências necessárias, incluindo alojamento e transporte."

# Paragraph 19: "We’re currently reviewing..."
Find-Replace $d.Paragraphs.Item(19).Range `
    "We’re currently reviewing your documents and will reach out to you if we need anything else. " `
    "Estamos a analisar os seus documentos e iremos contactá-lo caso seja necessária mais alguma informação. "

# Paragraph 20: "If you have any questions, please contact us via live chat or WhatsApp."
$para20 = $d.Paragraphs.Item(20).Range
Find-Replace $para20 `
    "If you have any questions, please contact us via " `
    "Para mais informações, contacte-nos através de "
Find-Replace $para20 " or " " ou "

# Paragraph 21: country manager contact details.
$para21 = $d.Paragraphs.Item(21).Range
Find-Replace $para21 `
    "If you have any questions, please contact your country manager, " `
    "Para mais questões, pode também contactar o seus gestor de parcerias "
Find-Replace $para21 ", at " ", em "
Find-Replace $para21 " or " " ou "

# Paragraph 22: "We look forward to seeing you at [EVENT NAME]. "
$para22 = $d.Paragraphs.Item(22).Range
Find-Replace $para22 `
    "We look forward to seeing you at " `
    "Esperamos vê-lo em breve, na "
Find-Replace $para22 "[EVENT NAME]" "[NOME DO EVENTO]"
